# Auto-generated Excel COM-interop edit script
# Applies per-cell value updates (and a few cell additions/removals)
# derived from the authoritative XML diff of Halicarnassus_Profits workbook sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 667.6667
$ws.Range("I6").Value = 434.83334
$ws.Range("K6").Value = 1304.50002
$ws.Range("M6").Value = -1192.50002
$ws.Range("H9").Value = 71.53846
$ws.Range("I9").Value = 60
$ws.Range("K9").Value = 60
$ws.Range("M9").Value = 109
$ws.Range("H12").Value = 200
$ws.Range("I12").Value = 200
$ws.Range("K12").Value = 200
$ws.Range("M12").Value = -30
$ws.Range("H21").Value = 1500
$ws.Range("I21").Value = 1500
$ws.Range("K21").Value = 1500
$ws.Range("M21").Value = -1032
$ws.Range("H23").Value = 1500
$ws.Range("I23").Value = 1500
$ws.Range("K23").Value = 1500
$ws.Range("M23").Value = -1266
$ws.Range("H28").Value = 4143.077
$ws.Range("I28").Value = 945.3333
$ws.Range("K28").Value = 945.3333
$ws.Range("M28").Value = -460.3333
$ws.Range("H29").Value = 4106.6523
$ws.Range("J29").Value = 7392
$ws.Range("L29").Value = 22176
$ws.Range("N29").Value = -22738
$ws.Range("H31").Value = 209.33333
$ws.Range("I31").Value = 209.33333
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 627.99999
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -397.99999
$ws.Range("N31").ClearContents()
$ws.Range("H38").Value = 4162.3125
$ws.Range("I38").Value = 3040.7
$ws.Range("K38").Value = 9122.099999999999
$ws.Range("M38").Value = -8750.099999999999
$ws.Range("H39").Value = 224.875
$ws.Range("J39").Value = 319.25
$ws.Range("L39").Value = 957.75
$ws.Range("N39").Value = -1549.75
$ws.Range("H40").Value = 4747.6665
$ws.Range("I40").Value = 3939.7058
$ws.Range("J40").Value = 7245
$ws.Range("K40").Value = 3939.7058
$ws.Range("L40").Value = 7245
$ws.Range("M40").Value = -3764.7058
$ws.Range("N40").Value = -7595
$ws.Range("H43").Value = 5364.8335
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H58").Value = 2700
$ws.Range("J58").Value = 2700
$ws.Range("L58").Value = 8100
$ws.Range("N58").Value = -8400
$ws.Range("H70").Value = 2900.2727
$ws.Range("J70").Value = 3569.8
$ws.Range("L70").Value = 10709.4
$ws.Range("N70").Value = -11249.4
$ws.Range("H73").Value = 2900.2727
$ws.Range("J73").Value = 3569.8
$ws.Range("L73").Value = 10709.4
$ws.Range("N73").Value = -12581.4
$ws.Range("H96").Value = 982.26666
$ws.Range("I96").Value = 654.625
$ws.Range("J96").Value = 1356.7142
$ws.Range("K96").Value = 1963.875
$ws.Range("L96").Value = 4070.1426
$ws.Range("M96").Value = -590.875
$ws.Range("N96").Value = -6816.142599999999
$ws.Range("H100").Value = 2149.1667
$ws.Range("I100").Value = 2149.1667
$ws.Range("K100").Value = 2149.1667
$ws.Range("M100").Value = -1608.1667
$ws.Range("H107").Value = 1043.85
$ws.Range("I107").Value = 1038.7778
$ws.Range("K107").Value = 1038.7778
$ws.Range("M107").Value = 881.2221999999999
$ws.Range("H111").Value = 686.6667
$ws.Range("I111").Value = 380
$ws.Range("J111").Value = 1300
$ws.Range("K111").Value = 1140
$ws.Range("L111").Value = 3900
$ws.Range("M111").Value = 1927
$ws.Range("N111").Value = -10034
$ws.Range("H116").Value = 2185.2856
$ws.Range("I116").Value = 2063.5454
$ws.Range("K116").Value = 2063.5454
$ws.Range("M116").Value = 1378.4546
$ws.Range("H125").Value = 5016.6665
$ws.Range("I125").Value = 4400
$ws.Range("J125").Value = 5325
$ws.Range("K125").Value = 39600
$ws.Range("L125").Value = 47925
$ws.Range("M125").Value = -37140
$ws.Range("N125").Value = -52845
$ws.Range("H140").Value = 49800
$ws.Range("J140").Value = 49800
$ws.Range("L140").Value = 49800
$ws.Range("N140").Value = -60160

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 91.75
$ws.Range("I4").Value = 85
$ws.Range("J4").Value = 98.5
$ws.Range("K4").Value = 85
$ws.Range("L4").Value = 98.5
$ws.Range("M4").Value = 31
$ws.Range("N4").Value = -330.5
$ws.Range("H11").Value = 2216.6
$ws.Range("I11").Value = 801.5
$ws.Range("K11").Value = 801.5
$ws.Range("M11").Value = -657.5
$ws.Range("H45").Value = 2629.6667
$ws.Range("I45").Value = 1729.8334
$ws.Range("J45").Value = 4429.3335
$ws.Range("K45").Value = 1729.8334
$ws.Range("L45").Value = 4429.3335
$ws.Range("M45").Value = -1352.8334
$ws.Range("N45").Value = -5183.3335
$ws.Range("H88").Value = 1460.7778
$ws.Range("J88").Value = 1424.6
$ws.Range("L88").Value = 1424.6
$ws.Range("N88").Value = -2236.6
$ws.Range("H91").Value = 1460.7778
$ws.Range("J91").Value = 1424.6
$ws.Range("L91").Value = 1424.6
$ws.Range("N91").Value = -4232.6
$ws.Range("H97").Value = 551.0625
$ws.Range("I97").Value = 768.8182
$ws.Range("K97").Value = 768.8182
$ws.Range("M97").Value = -272.8182
$ws.Range("H110").Value = 1030.5714
$ws.Range("I110").Value = 765.4545000000001
$ws.Range("J110").Value = 2002.6666
$ws.Range("K110").Value = 765.4545000000001
$ws.Range("L110").Value = 2002.6666
$ws.Range("M110").Value = 1279.5455
$ws.Range("N110").Value = -6092.6666
$ws.Range("H132").Value = 2535.1667
$ws.Range("I132").Value = 2535.1667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7605.500100000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5075.500100000001
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 5100
$ws.Range("I5").Value = 7525
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 7525
$ws.Range("L5").Value = 250
$ws.Range("M5").Value = -7412
$ws.Range("N5").Value = -476
$ws.Range("H11").Value = 314.75
$ws.Range("I11").Value = 314.75
$ws.Range("K11").Value = 314.75
$ws.Range("M11").Value = -174.75
$ws.Range("H20").Value = 2882.7144
$ws.Range("I20").Value = 2069
$ws.Range("K20").Value = 2069
$ws.Range("M20").Value = -1822
$ws.Range("H54").Value = 2483.1667
$ws.Range("I54").Value = 1779.8
$ws.Range("K54").Value = 1779.8
$ws.Range("M54").Value = -1295.8
$ws.Range("H80").Value = 1643.2
$ws.Range("I80").Value = 1172
$ws.Range("J80").Value = 2350
$ws.Range("K80").Value = 1172
$ws.Range("L80").Value = 2350
$ws.Range("M80").Value = -174
$ws.Range("N80").Value = -4346
$ws.Range("H82").Value = 27363.4
$ws.Range("I82").Value = 7957.875
$ws.Range("K82").Value = 7957.875
$ws.Range("M82").Value = -7574.875
$ws.Range("H83").Value = 1643.2
$ws.Range("I83").Value = 1172
$ws.Range("J83").Value = 2350
$ws.Range("K83").Value = 5860
$ws.Range("L83").Value = 11750
$ws.Range("M83").Value = -868
$ws.Range("N83").Value = -21734
$ws.Range("H85").Value = 27363.4
$ws.Range("I85").Value = 7957.875
$ws.Range("K85").Value = 7957.875
$ws.Range("M85").Value = -6631.875
$ws.Range("H86").Value = 5205.56
$ws.Range("I86").Value = 4480.846
$ws.Range("J86").Value = 5990.6665
$ws.Range("K86").Value = 4480.846
$ws.Range("L86").Value = 5990.6665
$ws.Range("M86").Value = -3357.846
$ws.Range("N86").Value = -8236.666499999999
$ws.Range("H89").Value = 5205.56
$ws.Range("I89").Value = 4480.846
$ws.Range("J89").Value = 5990.6665
$ws.Range("K89").Value = 22404.23
$ws.Range("L89").Value = 29953.3325
$ws.Range("M89").Value = -16788.23
$ws.Range("N89").Value = -41185.3325
$ws.Range("H97").Value = 11650
$ws.Range("I97").Value = 11650
$ws.Range("K97").Value = 11650
$ws.Range("M97").Value = -10659

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 6711.9375
$ws.Range("I7").Value = 7780.077
$ws.Range("J7").Value = 2083.3333
$ws.Range("K7").Value = 7780.077
$ws.Range("L7").Value = 2083.3333
$ws.Range("M7").Value = -7667.077
$ws.Range("N7").Value = -2309.3333
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()
$ws.Range("H58").Value = 3787.077
$ws.Range("J58").Value = 6972.5
$ws.Range("L58").Value = 6972.5
$ws.Range("N58").Value = -7378.5
$ws.Range("H99").Value = 2313.75
$ws.Range("I99").Value = 1986
$ws.Range("J99").Value = 2400
$ws.Range("K99").Value = 1986
$ws.Range("L99").Value = 2400
$ws.Range("M99").Value = -488
$ws.Range("N99").Value = -5396
$ws.Range("H105").Value = 2873.375
$ws.Range("I105").Value = 2797.6
$ws.Range("J105").Value = 2999.6667
$ws.Range("K105").Value = 2797.6
$ws.Range("L105").Value = 2999.6667
$ws.Range("M105").Value = -1050.6
$ws.Range("N105").Value = -6493.6667
$ws.Range("H122").Value = 2362.5
$ws.Range("I122").Value = 2362.5
$ws.Range("K122").Value = 7087.5
$ws.Range("M122").Value = -4637.5
$ws.Range("H126").Value = 2313.75
$ws.Range("I126").Value = 1986
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 5958
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -3488
$ws.Range("N126").Value = -12140
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -470
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 1854.1333
$ws.Range("I134").Value = 1854.1333
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5562.3999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3027.3999
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 3787.077
$ws.Range("J136").Value = 6972.5
$ws.Range("L136").Value = 20917.5
$ws.Range("N136").Value = -26017.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 31.833334
$ws.Range("I6").Value = 17.75
$ws.Range("K6").Value = 53.25
$ws.Range("M6").Value = 59.75
$ws.Range("H7").Value = 51.25
$ws.Range("J7").Value = 56.666668
$ws.Range("L7").Value = 170.000004
$ws.Range("N7").Value = -394.000004
$ws.Range("H12").Value = 100.42105
$ws.Range("I12").Value = 14.75
$ws.Range("K12").Value = 44.25
$ws.Range("M12").Value = 128.75
$ws.Range("H38").Value = 297.51852
$ws.Range("I38").Value = 320.85715
$ws.Range("J38").Value = 215.83333
$ws.Range("K38").Value = 962.5714499999999
$ws.Range("L38").Value = 647.49999
$ws.Range("M38").Value = -615.5714499999999
$ws.Range("N38").Value = -1341.49999
$ws.Range("H44").Value = 20
$ws.Range("I44").Value = 20
$ws.Range("K44").Value = 60
$ws.Range("M44").Value = 338
$ws.Range("H80").Value = 4261.1924
$ws.Range("J80").Value = 4643.4443
$ws.Range("L80").Value = 13930.3329
$ws.Range("N80").Value = -15802.3329
$ws.Range("H83").Value = 4261.1924
$ws.Range("J83").Value = 4643.4443
$ws.Range("L83").Value = 41790.9987
$ws.Range("N83").Value = -51150.9987
$ws.Range("H92").Value = 864.25
$ws.Range("J92").Value = 977.5
$ws.Range("L92").Value = 2932.5
$ws.Range("N92").Value = -5428.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 19866.666
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 19866.666
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 19866.666
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -20168.666
$ws.Range("H70").Value = 1500
$ws.Range("I70").Value = 1500
$ws.Range("K70").Value = 1500
$ws.Range("M70").Value = -1230
$ws.Range("H73").Value = 1500
$ws.Range("I73").Value = 1500
$ws.Range("K73").Value = 1500
$ws.Range("M73").Value = -564
$ws.Range("H102").Value = 1709.3158
$ws.Range("I102").Value = 1322.2354
$ws.Range("K102").Value = 1322.2354
$ws.Range("M102").Value = 299.7646
$ws.Range("H126").Value = 2955
$ws.Range("I126").Value = 2955
$ws.Range("K126").Value = 8865
$ws.Range("M126").Value = -6395

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 39888.668
$ws.Range("J2").Value = 72999.5
$ws.Range("L2").Value = 72999.5
$ws.Range("N2").Value = -73223.5
$ws.Range("H22").Value = 2000
$ws.Range("J22").Value = 2000
$ws.Range("L22").Value = 2000
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 2000
$ws.Range("J27").Value = 2000
$ws.Range("L27").Value = 2000
$ws.Range("N27").Value = -2214
$ws.Range("H93").Value = 564.6667
$ws.Range("I93").Value = 564.6667
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 564.6667
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 683.3333
$ws.Range("N93").ClearContents()
$ws.Range("H132").Value = 3975
$ws.Range("I132").Value = 3975
$ws.Range("K132").Value = 11925
$ws.Range("M132").Value = -9395
$ws.Range("H136").Value = 3031
$ws.Range("I136").Value = 3031
$ws.Range("K136").Value = 9093
$ws.Range("M136").Value = -6543

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H81").Value = 6125
$ws.Range("J81").Value = 10000
$ws.Range("L81").Value = 20000
$ws.Range("N81").Value = -22122
$ws.Range("H84").Value = 6125
$ws.Range("J84").Value = 10000
$ws.Range("L84").Value = 100000
$ws.Range("N84").Value = -110608
$ws.Range("H96").Value = 200
$ws.Range("J96").Value = 200
$ws.Range("L96").Value = 200
$ws.Range("N96").Value = -2946
$ws.Range("H113").Value = 595
$ws.Range("I113").Value = 576
$ws.Range("K113").Value = 1728
$ws.Range("M113").Value = 442
$ws.Range("H126").Value = 5213.684
$ws.Range("I126").Value = 2706
$ws.Range("K126").Value = 8118
$ws.Range("M126").Value = -5648
$ws.Range("H132").Value = 1894.6471
$ws.Range("I132").Value = 1894.6471
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5683.9413
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3153.9413
$ws.Range("N132").ClearContents()
